$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: updated timestamp text ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 11:52"

# --- Row 15: Belgica - refreshed daily figures ---
$ws.Range("B15").Value = 41889
$ws.Range("C15").Value = 933
$ws.Range("D15").Value = 9433
$ws.Range("E15").Value = 26194
$ws.Range("F15").Value = 1020
$ws.Range("G15").Value = 264
$ws.Range("H15").Value = 6262

# --- Rows 62/63: Kuwait overtakes Hungria in ranking, so they swap rows ---
# Row 62 becomes Kuwait with refreshed figures
$ws.Range("A62").Value = "Kuwait"
$ws.Range("B62").Value = 2248
$ws.Range("C62").Value = 168
$ws.Range("D62").Value = 443
$ws.Range("E62").Value = 1792
$ws.Range("F62").Value = 50
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 13

# Row 63 becomes Hungria, carrying its previous figures
$ws.Range("A63").Value = "Hungria"
$ws.Range("B63").Value = 2168
$ws.Range("C63").Value = 70
$ws.Range("D63").Value = 295
$ws.Range("E63").Value = 1648
$ws.Range("F63").Value = 82
$ws.Range("G63").Value = 12
$ws.Range("H63").Value = 225

# --- Row 76: Eslovenia - refreshed daily figures ---
$ws.Range("B76").Value = 1353
$ws.Range("C76").Value = 9
$ws.Range("D76").Value = 205
$ws.Range("E76").Value = 1069
$ws.Range("F76").Value = 24
$ws.Range("G76").Value = 2
$ws.Range("H76").Value = 79

# --- Row 95: Libano - refreshed daily figures ---
$ws.Range("B95").Value = 682
$ws.Range("C95").Value = 5
$ws.Range("D95").Value = 108
$ws.Range("E95").Value = 552
$ws.Range("F95").Value = 26
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 22

# --- Row 205: Papua Nueva Guinea - refreshed daily figures ---
$ws.Range("B205").Value = 8
$ws.Range("C205").Value = 1
$ws.Range("D205").Value = 0
$ws.Range("E205").Value = 8
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0
